# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets
# to the newly scraped numbers.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for sheet "展览"
$exhibitUpdates = @{
    6  = 2968
    8  = 1979
    10 = 297
    11 = 809
    12 = 945
    13 = 197
    14 = 412
    15 = 1134
    19 = 7154
    21 = 1844
    22 = 187
    24 = 162
    25 = 382
    26 = 356
    28 = 1118
    29 = 941
    31 = 116
    33 = 1111
    34 = 1918
    35 = 469
    36 = 7
    38 = 250
    41 = 268
}

# Row -> new F-column value for sheet "全部类型"
$allTypesUpdates = @{
    9  = 2968
    11 = 1979
    13 = 297
    14 = 809
    16 = 945
    17 = 197
    18 = 412
    19 = 1134
    23 = 7154
    25 = 1844
    27 = 187
    29 = 162
    30 = 382
    31 = 356
    33 = 1118
    34 = 941
    36 = 116
    37 = 1111
    38 = 1918
    39 = 469
    40 = 7
    42 = 250
    45 = 268
}

$wsExhibit = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
